$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exp_meta")

$ws.Range('CR2').Value = 'Q_EUROSTAT'
$ws.Range('CS2').Value = 'Q_DESTATIS'
$ws.Range('CR3').Value = 'L_BMEL_1'
$ws.Range('CS3').Value = 'Q_JKI'
$ws.Range('CT3').Value = 'Q_UG'
$ws.Range('CR4').Value = 'L_ERSTT_1'
$ws.Range('CS4').Value = 'L_BMEL_2'
$ws.Range('CT4').Value = 'L_DSTTS_7'
$ws.Range('CR5').Value = 'Q_BMZ'
$ws.Range('CS5').Value = 'Q_BMEL'
$ws.Range('CR6').Value = 'L_GBE_1'
$ws.Range('CR7').Value = 'L_BZGA_1'
$ws.Range('CS7').Value = 'Q_DESTATIS'
$ws.Range('CR8').Value = 'Q_RKI'
$ws.Range('CR9').Value = 'Q_DESTATIS'
$ws.Range('CR10').Value = 'L_UBA_1'
$ws.Range('CR11').Value = 'Q_UBA'
$ws.Range('CS11').Value = 'Q_WHO'
$ws.Range('CR12').Value = 'Q_BMG'
$ws.Range('CS12').Value = 'Q_AA'
$ws.Range('CT12').Value = 'Q_BMZ'
$ws.Range('CU12').Value = 'Q_BMBF'
$ws.Range('CR13').Value = 'Q_DESTATIS'
$ws.Range('CR14').Value = 'Q_DESTATIS'
$ws.Range('CR15').Value = 'Q_DESTATIS'
$ws.Range('CR16').Value = 'Q_BMZ'
$ws.Range('CS16').Value = 'Q_CEVAL'
$ws.Range('CR17').Value = 'L_DSTTS_1'
$ws.Range('CR18').Value = 'L_FIDAR_1'
$ws.Range('CS18').Value = 'Q_DESTATIS'
$ws.Range('CR19').Value = 'L_DSTTS_2'
$ws.Range('CR20').Value = 'L_UBA_2'
$ws.Range('CR21').Value = 'L_UBA_3'
$ws.Range('CS21').Value = 'L_LANUV_1'
$ws.Range('CR22').Value = 'Q_KFW'
$ws.Range('CR26').Value = 'Q_DESTATIS'
$ws.Range('CR27').Value = 'L_DSTTS_3'
$ws.Range('CS27').Value = 'Q_BMF'
$ws.Range('CR28').Value = 'L_BB_1'
$ws.Range('CS28').Value = 'Q_DESTATIS'
$ws.Range('CR29').Value = 'Q_DESTATIS'
$ws.Range('CS29').Value = 'L_SP_1'
$ws.Range('CR30').Value = 'Q_DESTATIS'
$ws.Range('CS30').Value = 'L_SP_2'
$ws.Range('CR31').Value = 'Q_DESTATIS'
$ws.Range('CS31').Value = 'L_ERSTT_2'
$ws.Range('CT31').Value = 'Q_AGNE'
$ws.Range('CR32').Value = 'Q_GIZ'
$ws.Range('CR33').Value = 'L_DSTTS_4'
$ws.Range('CS33').Value = 'L_DSTTS_9'
$ws.Range('CT33').Value = 'Q_VWGDL'
$ws.Range('CU33').Value = 'Q_SVWS'
$ws.Range('CR34').Value = 'Q_BMDV'
$ws.Range('CR35').Value = 'Q_DESTATIS'
$ws.Range('CR36').Value = 'L_DSTTS_5'
$ws.Range('CS36').Value = 'L_ERSTT_3'
$ws.Range('CT36').Value = 'L_ERSTT_6'
$ws.Range('CU36').Value = 'L_BB_2'
$ws.Range('CV36').Value = 'L_EZB_1'
$ws.Range('CW36').Value = 'Q_DIW'
$ws.Range('CR37').Value = 'Q_DESTATIS'
$ws.Range('CS37').Value = 'L_LANUV_2'
$ws.Range('CR38').Value = 'Q_DESTATIS'
$ws.Range('CS38').Value = 'Q_BBSR'
$ws.Range('CT38').Value = 'Q_JHT'
$ws.Range('CR39').Value = 'Q_DESTATIS'
$ws.Range('CS39').Value = 'Q_BBSR'
$ws.Range('CT39').Value = 'Q_JHT'
$ws.Range('CR40').Value = 'Q_DESTATIS'
$ws.Range('CS40').Value = 'Q_IFEU'
$ws.Range('CR41').Value = 'Q_DESTATIS'
$ws.Range('CS41').Value = 'Q_IFEU'
$ws.Range('CR42').Value = 'Q_BBSR'
$ws.Range('CS42').Value = 'Q_IFEU'
$ws.Range('CR43').Value = 'L_DSTTS_6'
$ws.Range('CS43').Value = 'L_ERSTT_4'
$ws.Range('CR44').Value = 'Q_BKM'
$ws.Range('CR45').Value = 'Q_DESTATIS'
$ws.Range('CR46').Value = 'Q_GFK'
$ws.Range('CS46').Value = 'Q_KBA'
$ws.Range('CT46').Value = 'Q_AMI'
$ws.Range('CU46').Value = 'Q_BOLW'
$ws.Range('CV46').Value = 'Q_VCD'
$ws.Range('CW46').Value = 'L_UBA_4'
$ws.Range('CR47').Value = 'L_EMAS_1'
$ws.Range('CR48').Value = 'Q_KNB'
$ws.Range('CS48').Value = 'Q_IFEU'
$ws.Range('CT48').Value = 'Q_UBA'
$ws.Range('CU48').Value = 'Q_DESTATIS'
$ws.Range('CR49').Value = 'L_UBA_5'
$ws.Range('CS49').Value = 'L_SP_3'
$ws.Range('CR50').Value = 'Q_BMZ'
$ws.Range('CR51').Value = 'L_UBA_6'
$ws.Range('CR52').Value = 'L_ERSTT_5'
$ws.Range('CR53').Value = 'Q_BFN'
$ws.Range('CR54').Value = 'L_UBA_7'
$ws.Range('CR55').Value = 'Q_DESTATIS'
$ws.Range('CS55').Value = 'Q_BMZ'
$ws.Range('CR56').Value = 'L_BKA_1'
$ws.Range('CS56').Value = 'Q_DESTATIS'
$ws.Range('CR57').Value = 'Q_AA'
$ws.Range('CR58').Value = 'L_TI_1'
$ws.Range('CS58').Value = 'Q_BMZ'
$ws.Range('CR59').Value = 'Q_DESTATIS'
$ws.Range('CS59').Value = 'Q_BMZ'
$ws.Range('CT59').Value = 'Q_OECD'
$ws.Range('CR60').Value = 'Q_DESTATIS'
$ws.Range('CR61').Value = 'Q_DESTATIS'
